$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 226
$wsExhibit.Range("F4").Value = 2544
$wsExhibit.Range("F6").Value = 548

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 226
$wsAll.Range("F6").Value = 2544
$wsAll.Range("F8").Value = 548
